$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1, reusing the exact formatting (style) of the
# existing header row (copy format only from O1, so no new style entries are
# introduced in styles.xml).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# Update existing columns I, K, M, O and populate new columns P, Q for rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I - was 1
    $ws.Cells.Item($r, 11).Value = 1   # K - was 2
    $ws.Cells.Item($r, 13).Value = 2   # M - was 1
    $ws.Cells.Item($r, 15).Value = 1   # O - was 2
    $ws.Cells.Item($r, 16).Value = 2   # P - new
    $ws.Cells.Item($r, 17).Value = 2   # Q - new
}
